# Update the footer placeholder on every slide:
#   "107th IETF @ Vancouver"  ->  "108th IETF @ Madrid"
# Runs are edited in-place via TextRange.Characters(start, length) so the
# existing run-level formatting (e.g. the superscript "th" run, sz="1200"
# runs, etc.) is preserved exactly as in the original file.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)

        if ($sh.Name -ne "Footer Placeholder 3") {
            continue
        }
        if (-not $sh.HasTextFrame) {
            continue
        }

        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text

        # Replace the "107" year number, wherever it sits in the string.
        $numIdx = $full.IndexOf("107")
        if ($numIdx -ge 0) {
            $numRange = $tr.Characters($numIdx + 1, 3)
            $numRange.Text = "108"
        }

        # Re-read the text (it changed above) before locating the city.
        $full = $tr.Text
        $cityIdx = $full.IndexOf(" IETF @ Vancouver")
        if ($cityIdx -ge 0) {
            $cityRange = $tr.Characters($cityIdx + 1, " IETF @ Vancouver".Length)
            $cityRange.Text = " IETF @ Madrid"
        }
    }
}
